$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SSN rows (row 10): data type column changes from "int" to "Varchar(9)"
# for both the Users table (B10) and the Applications table (F10).
$ws.Range("B10").Value = "Varchar(9)"
$ws.Range("F10").Value = "Varchar(9)"

# New rows appended to the "Addresses" table (columns A-C, below the
# existing acc_num / appli_num rows at 17-18).
$ws.Range("A19").Value = "street add"
$ws.Range("B19").Value = "Varchar(100)"
$ws.Range("C19").Value = "NOT NULL"

$ws.Range("A20").Value = "city"
$ws.Range("B20").Value = "Varchar(100)"
$ws.Range("C20").Value = "NOT NULL"

$ws.Range("A21").Value = "state"
$ws.Range("B21").Value = "Varchar(100)"
$ws.Range("C21").Value = "NOT NULL"

$ws.Range("A22").Value = "zip_code"
$ws.Range("B22").Value = "int"
$ws.Range("C22").Value = "NOT NULL"

$ws.Range("A23").Value = "country"
$ws.Range("B23").Value = "Varchar(100) "
$ws.Range("C23").Value = "NOT NULL"

# Move the active selection to A24, matching the post-edit cursor position.
$ws.Range("A24").Select()
